$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'310.82"
$ws.Range("E2").Value = "'0.66%"
$ws.Range("D3").Value = "'39.91"
$ws.Range("E3").Value = "'-2.85%"
$ws.Range("D4").Value = "'5.102"
$ws.Range("E4").Value = "'-2.54%"
$ws.Range("E5").Value = "'-1.23%"
$ws.Range("D6").Value = "'4.300"
$ws.Range("E6").Value = "'-0.47%"
$ws.Range("D7").Value = "'1.663"
$ws.Range("E7").Value = "'2.62%"
$ws.Range("D8").Value = "'0.9295"
$ws.Range("E8").Value = "'1.20%"
$ws.Range("E9").Value = "'-0.86%"
$ws.Range("D10").Value = "'0.1214"
$ws.Range("E10").Value = "'-2.00%"
$ws.Range("D11").Value = "'0.1814"
$ws.Range("E11").Value = "'-1.35%"
$ws.Range("D12").Value = "'0.09021"
$ws.Range("E12").Value = "'-1.68%"
$ws.Range("D13").Value = "'0.04152"
$ws.Range("E13").Value = "'-3.23%"
$ws.Range("E14").Value = "'0.17%"
$ws.Range("D15").Value = "'0.001284"
$ws.Range("E15").Value = "'1.74%"
$ws.Range("D16").Value = "'0.005961"
$ws.Range("E16").Value = "'2.48%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.351"
$ws.Range("E18").Value = "'-0.04%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3355"
$ws.Range("E19").Value = "'0.58%"
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").Value = "'7.658"
$ws.Range("E20").Value = "'6.13%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1352"
$ws.Range("E21").Value = "'-2.28%"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "'0.2811"
$ws.Range("E22").Value = "'-2.84%"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").Value = "'0.04023"
$ws.Range("E23").Value = "'-1.34%"
$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D24").Value = "'0.001267"
$ws.Range("E24").Value = "'0.39%"
$ws.Range("B25").Value = "HotbitToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D25").Value = "'0.004098"
$ws.Range("E25").Value = "'-1.01%"
$ws.Range("D26").Value = "'0.0001271"
$ws.Range("E26").Value = "'-0.11%"
$ws.Range("D38").Value = "'0.02426"
$ws.Range("E38").Value = "'-0.95%"
$ws.Range("D39").Value = "'0.05142"
$ws.Range("E39").Value = "'-2.97%"
$ws.Range("D40").Value = "'0.007734"
$ws.Range("E40").Value = "'-1.51%"
$ws.Range("D41").Value = "'0.1298"
$ws.Range("E41").Value = "'-1.13%"
$ws.Range("D42").Value = "'0.007659"
$ws.Range("E42").Value = "'12.26%"
$ws.Range("E43").Value = "'14.25%"
$ws.Range("D44").Value = "'0.008070"
$ws.Range("E44").Value = "'3.60%"
$ws.Range("E45").Value = "'1.34%"
$ws.Range("D46").Value = "'0.00006588"
$ws.Range("E46").Value = "'-1.18%"
$ws.Range("E47").Value = "'-0.11%"
$ws.Range("D48").Value = "'0.2600"
$ws.Range("E48").Value = "'53.05%"
$ws.Range("D49").Value = "'0.004204"
$ws.Range("E49").Value = "'2.54%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'-0.11%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'-0.11%"
Write-Host "Updated symbol list"
